$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Higher ed must get ahead of AI paradigm shift (opinion)"
$ws.Range("B2").Value = "* AI is rapidly changing the definition of knowledge work and potentially displacing human workers.`n* Higher education needs to address the larger questions posed by AI, not just academic integrity.`n* Competency-based education models may become more relevant as AI handles basic cognitive tasks.`n* Curricula need to be updated to reflect the changing demands of the workforce in the age of AI.`n* Universities can leverage AI for administrative processes, instruction, and tutoring, but must also address potential job displacement.`n* Society needs to value and invest in distinctly human jobs that algorithms cannot do, such as teaching, social work, and mental health care.`n* Universities need philosophers, ethicists, historians, sociologists, psychologists, cognitive scientists, and law faculty to address the complexities of AI.`n* Interdisciplinary collaboration is crucial to understanding and managing the impact of AI."
$ws.Range("C2").Value = "* Addresses leveraging AI in education, specifically how AI tools are changing higher education and potentially K-12.`n* Discusses career readiness as AI changes the skills needed for the workforce.`n* Touches on the potential displacement of workers, which is relevant to career readiness and soft skill development.`n* Highlights the need for uniquely human skills, implicitly emphasizing soft skill development.`n* Examines the role of edtechs in instruction, tutoring and advising.`n* Raises important questions regarding learning design and how AI may change the role of knowledge in learning."
